# Add a new student record (row 11) to the registration form demo sheet.
# Values are entered in the same order the original author typed the
# registration form fields, so new shared-string entries land in the
# same sequence as the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# stu_first_name
$ws.Range("A11").Value = "Pulkit"

# stu_father_name
$ws.Range("F11").Value = "Rahul Aggarwal"

# stu_last_name
$ws.Range("B11").Value = "Aggarwal"

# stu_mother_name
$ws.Range("H11").Value = "Rajni Rahul Aggarwal"

# stu_address_line1
$ws.Range("I11").Value = "989 Bazar Sita Ram"

# stu_address_line2
$ws.Range("J11").Value = "Chawri Bazar"

# stu_class
$ws.Range("D11").Value = "X"

# stu_roll_no
$ws.Range("E11").Value = 123354

# stu_contact
$ws.Range("G11").Value = 9213456859

# stu_city
$ws.Range("K11").Value = "Delhi"

# stu_postalcode
$ws.Range("L11").Value = 110006

# stu_full_name - same formula pattern used by the rows above it
$ws.Range("C11").Formula = "=A11&"" ""&B11"

# Match the author's final cursor position/selection on the sheet.
$ws.Range("I11").Select()
